# Fruta / hortaliza, semanal
# Re-orders the weekly price rows (2-17) by directly writing the
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Precio $/Kg) tuples to their new row positions. Row 4 and row 18 are
# unchanged by this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, columns: D, J, K, L, M, P
$rows = @{
    2  = @(44895, 30, 18000, 18000, 18000, 1385)
    3  = @(44880, 30, 17000, 17000, 17000, 1308)
    5  = @(44839, 40, 15000, 16000, 15500, 1192)
    6  = @(44797, 60, 12000, 13000, 12500, 962)
    7  = @(44943, 30, 17000, 17000, 17000, 1308)
    8  = @(44859, 30, 13000, 13000, 13000, 1000)
    9  = @(44841, 30, 18000, 18000, 18000, 1385)
    10 = @(44874, 30, 17000, 17000, 17000, 1308)
    11 = @(44868, 30, 18000, 18000, 18000, 1385)
    12 = @(44922, 30, 17000, 17000, 17000, 1308)
    13 = @(44915, 50, 18000, 18000, 18000, 1385)
    14 = @(44804, 40, 12000, 13000, 12500, 962)
    15 = @(44894, 30, 18000, 18000, 18000, 1385)
    16 = @(44846, 30, 18000, 18000, 18000, 1385)
    17 = @(44832, 60, 17000, 18000, 17500, 1346)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
